# Fruta / hortaliza, semanal
# Inserts a new weekly price record for Cilantro at Feria Lagunitas de Puerto Montt,
# shifting existing rows 268-356 down to 269-357.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 268, pushing the remaining data down.
$ws.Rows.Item(268).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A268").Value2() = 4
$ws.Range("B268").Value2() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C268").Value2() = "Los Lagos"
$ws.Range("D268").Value2() = 44876
$ws.Range("E268").Value2() = 10
$ws.Range("F268").Value2() = 100112040
$ws.Range("G268").Value2() = "Cilantro"
$ws.Range("H268").Value2() = "Sin especificar"
$ws.Range("I268").Value2() = "Primera"
$ws.Range("J268").Value2() = 160
$ws.Range("K268").Value2() = 8000
$ws.Range("L268").Value2() = 8000
$ws.Range("M268").Value2() = 8000
$ws.Range("N268").Value2() = "$/docena de atados (2 kilos)"
$ws.Range("O268").Value2() = "Región de La Araucanía"
$ws.Range("P268").Value2() = 4000
$ws.Range("Q268").Value2() = 2
$ws.Range("R268").Value2() = "Hortaliza"
